# This script applies a row-content rotation/swap within the Artfynd export sheet.
# Rows 5, 8, 9 undergo a 3-way rotation of all their data (row 5 <- old row 8,
# row 8 <- old row 9, row 9 <- old row 5). Rows 12/13 and rows 16/17 each swap
# their full content pairwise. The row position (r="N") itself never moves;
# only the cell contents within each row are reassigned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 5 ----
$ws.Range("A5").Value = 130822173
$ws.Range("B5").Value = 91828
$ws.Range("E5").Value = 5432
$ws.Range("F5").Value = 'Granticka'
$ws.Range("G5").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H5").Value = ""
$ws.Range("K5").Value = 'teleomorf'
$ws.Range("Q5").Value = 426315
$ws.Range("R5").Value = 7048136
$ws.Range("AC5").Value = 'Flera fruktkroppar intill andra vedsvampar i en ca 5 meters granhögstubbe.'
$ws.Range("AM5").Value = 'Stående död trädstam/högstubbe'
$ws.Range("AO5").Value = 'Standing dead tree/snags # Picea abies'

# ---- Row 8 ----
$ws.Range("A8").Value = 130822207
$ws.Range("B8").Value = 79243
$ws.Range("E8").Value = 6425
$ws.Range("F8").Value = 'Garnlav'
$ws.Range("G8").Value = 'Alectoria sarmentosa'
$ws.Range("H8").Value = '(Ach.) Ach.'
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = ""
$ws.Range("N8").Value = ""
$ws.Range("P8").Value = 'Djupsjön ö, Jmt'
$ws.Range("Q8").Value = 426504
$ws.Range("R8").Value = 7047832
$ws.Range("AC8").Value = 'Måttligt'
$ws.Range("AF8").Value = ""
$ws.Range("AH8").Value = ""
$ws.Range("AJ8").Value = ""
$ws.Range("AK8").Value = ""
$ws.Range("AM8").Value = ""
$ws.Range("AO8").Value = ""
$ws.Range("AW8").Value = 'Benny Öwre'
$ws.Range("AX8").Value = 'Benny Öwre'

# ---- Row 9 ----
$ws.Range("A9").Value = 130822198
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = ""
$ws.Range("N9").Value = ""
$ws.Range("P9").Value = 'Djupsjön Öst, Jmt'
$ws.Range("Q9").Value = 426516
$ws.Range("R9").Value = 7048217
$ws.Range("AC9").Value = ""
$ws.Range("AF9").Value = ""
$ws.Range("AH9").Value = 'Granskog'
$ws.Range("AJ9").Value = 'gran'
$ws.Range("AK9").Value = 'Picea abies'
$ws.Range("AM9").Value = 'Gren på levande träd'
$ws.Range("AO9").Value = 'Branch on living tree # Picea abies'
$ws.Range("AW9").Value = 'Kristian Zackrisson'
$ws.Range("AX9").Value = 'Kristian Zackrisson'

# ---- Row 12 ----
$ws.Range("A12").Value = 130822177
$ws.Range("B12").Value = 79243
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = 'Garnlav'
$ws.Range("G12").Value = 'Alectoria sarmentosa'
$ws.Range("H12").Value = '(Ach.) Ach.'
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("P12").Value = 'Djupsjön Öst, Jmt'
$ws.Range("Q12").Value = 426500
$ws.Range("R12").Value = 7048040
$ws.Range("AC12").Value = 'På flera granar.'
$ws.Range("AF12").Value = ""
$ws.Range("AH12").Value = 'Granskog'
$ws.Range("AJ12").Value = 'gran'
$ws.Range("AK12").Value = 'Picea abies'
$ws.Range("AM12").Value = 'Gren på levande träd'
$ws.Range("AO12").Value = 'Branch on living tree # Picea abies'
$ws.Range("AW12").Value = 'Kristian Zackrisson'
$ws.Range("AX12").Value = 'Kristian Zackrisson'

# ---- Row 13 ----
$ws.Range("A13").Value = 130822204
$ws.Range("B13").Value = 91828
$ws.Range("E13").Value = 5432
$ws.Range("F13").Value = 'Granticka'
$ws.Range("G13").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H13").Value = ""
$ws.Range("J13").Value = ""
$ws.Range("K13").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("P13").Value = 'Djupsjön ö, Jmt'
$ws.Range("Q13").Value = 426409
$ws.Range("R13").Value = 7048165
$ws.Range("AC13").Value = ""
$ws.Range("AF13").Value = ""
$ws.Range("AH13").Value = ""
$ws.Range("AJ13").Value = ""
$ws.Range("AK13").Value = ""
$ws.Range("AM13").Value = ""
$ws.Range("AO13").Value = ""
$ws.Range("AW13").Value = 'Benny Öwre'
$ws.Range("AX13").Value = 'Benny Öwre'

# ---- Row 16 ----
$ws.Range("A16").Value = 130822203
$ws.Range("B16").Value = 91828
$ws.Range("E16").Value = 5432
$ws.Range("F16").Value = 'Granticka'
$ws.Range("G16").Value = 'Porodaedalea chrysoloma s.lat.'
$ws.Range("H16").Value = ""
$ws.Range("P16").Value = 'Djupsjön ö, Jmt'
$ws.Range("Q16").Value = 426286
$ws.Range("R16").Value = 7048138
$ws.Range("AC16").Value = ""
$ws.Range("AW16").Value = 'Benny Öwre'
$ws.Range("AX16").Value = 'Benny Öwre'

# ---- Row 17 ----
$ws.Range("A17").Value = 130822175
$ws.Range("B17").Value = 79243
$ws.Range("E17").Value = 6425
$ws.Range("F17").Value = 'Garnlav'
$ws.Range("G17").Value = 'Alectoria sarmentosa'
$ws.Range("H17").Value = '(Ach.) Ach.'
$ws.Range("P17").Value = 'Djupsjön Öst, Jmt'
$ws.Range("Q17").Value = 426484
$ws.Range("R17").Value = 7048011
$ws.Range("AC17").Value = 'Relativt rikligt på flera granar.'
$ws.Range("AW17").Value = 'Kristian Zackrisson'
$ws.Range("AX17").Value = 'Kristian Zackrisson'

Write-Output "Row rotation/swap applied for rows 5,8,9,12,13,16,17."
